$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# Update the daily conversion note on Hoja1!A1
$oldText = $ws1.Range("A1").Value2
$newText = $oldText.Replace("1000 Bs = 2.94 = 10895.85 pesos", "1000 Bs = 3.0 = 11160.15 pesos")
$newText = $newText.Replace("10895.85 pesos = 2.92 = 934.14 Bs", "11160.15 pesos = 2.99 = 956.42 Bs")
$ws1.Range("A1").Value = $newText

# Update the rate figures on the tasas sheet
$ws2.Range("N10").Value = 333.777
$ws2.Range("O10").Value = 3725
$ws2.Range("N12").Value = 3733.99
